$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.228
$ws.Range("C6").Value = -11.864
$ws.Range("C7").Value = -12.793
$ws.Range("D7").Value = -7.518000000000001
$ws.Range("D12").Value = -7.470999999999999
$ws.Range("E13").Value = 16.354
$ws.Range("E14").Value = 16.854
$ws.Range("D15").Value = -8.422999999999998
$ws.Range("C16").Value = -12.808
$ws.Range("E16").Value = 16.729
$ws.Range("E19").Value = 16.539
$ws.Range("C20").Value = -12.061
$ws.Range("D20").Value = -7.994000000000002
$ws.Range("D21").Value = -7.967000000000001
$ws.Range("D22").Value = -7.906000000000001
$ws.Range("E22").Value = 16.625
$ws.Range("D23").Value = -7.997
$ws.Range("C28").Value = -12.848
$ws.Range("C29").Value = -11.905
$ws.Range("D29").Value = -7.325
$ws.Range("C32").Value = -12.959
$ws.Range("D34").Value = -7.937
$ws.Range("E36").Value = 16.804
$ws.Range("C40").Value = -12.2
$ws.Range("D42").Value = -8.18
$ws.Range("D43").Value = -7.795
$ws.Range("D44").Value = -7.580999999999999
$ws.Range("D45").Value = -7.704000000000001
$ws.Range("C46").Value = -13.322
$ws.Range("D46").Value = -8.140000000000001
$ws.Range("E46").Value = 16.889
$ws.Range("D50").Value = -8.217000000000002
$ws.Range("E50").Value = 16.688
$ws.Range("C51").Value = -12.199
$ws.Range("D51").Value = -7.598999999999999
$ws.Range("C52").Value = -11.614
$ws.Range("C57").Value = -13.844
$ws.Range("C59").Value = -12.235
$ws.Range("C62").Value = -13.737
$ws.Range("C66").Value = -11.203
$ws.Range("D66").Value = -7.683000000000002
$ws.Range("D67").Value = -7.100999999999999
$ws.Range("C73").Value = -12.401
$ws.Range("C74").Value = -11.866
$ws.Range("D79").Value = -7.728
$ws.Range("D84").Value = -8.266000000000002
$ws.Range("C92").Value = -11.714
$ws.Range("D92").Value = -6.728
$ws.Range("E95").Value = 17.52
$ws.Range("D97").Value = -8.191000000000001
$ws.Range("E97").Value = 16.468
$ws.Range("C100").Value = -12.777
